# Egg laying datasheet finished.
# Append the 10.18.21 collection-date rows (167-177) to the
# mother_laying_bydate sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- collect_date column (B) -------------------------------------------
# Typing "10.18.21" straight into a cell gets auto-recognised as a date
# (like any d.mm.yy-ish literal), which is not how the existing rows are
# stored (they are plain shared-string text). To land it as literal text
# without flipping the cell's number format, stage it once as a text
# formula result in a scratch cell, then Paste-Special (values only) into
# every date cell that needs it - that copies the literal text through,
# cell by cell, as a shared string.
$helper = $ws.Range("H1")
$helper.Formula = '="10.18.21"'
$helper.Copy()

$ws.Cells.Item(167, 2).PasteSpecial(-4163)
$ws.Cells.Item(168, 2).PasteSpecial(-4163)
$ws.Cells.Item(169, 2).PasteSpecial(-4163)
$ws.Cells.Item(170, 2).PasteSpecial(-4163)
$ws.Cells.Item(171, 2).PasteSpecial(-4163)
$ws.Cells.Item(172, 2).PasteSpecial(-4163)
$ws.Cells.Item(173, 2).PasteSpecial(-4163)
$ws.Cells.Item(174, 2).PasteSpecial(-4163)
$ws.Cells.Item(175, 2).PasteSpecial(-4163)
$ws.Cells.Item(176, 2).PasteSpecial(-4163)
$ws.Cells.Item(177, 2).PasteSpecial(-4163)

$helper.ClearContents()

# --- remaining columns: MID (A), n_eggs (C), n_viable (D), pop (E) -----
$ws.Cells.Item(167, 1).Value = 114
$ws.Cells.Item(167, 3).Value = 18
$ws.Cells.Item(167, 4).Value = 4
$ws.Cells.Item(167, 5).Value = "KL"

$ws.Cells.Item(168, 1).Value = 118
$ws.Cells.Item(168, 3).Value = 21
$ws.Cells.Item(168, 4).Value = 0
$ws.Cells.Item(168, 5).Value = "KL"

$ws.Cells.Item(169, 1).Value = 329
$ws.Cells.Item(169, 3).Value = 24
$ws.Cells.Item(169, 4).Value = 0
$ws.Cells.Item(169, 5).Value = "PK"

$ws.Cells.Item(170, 1).Value = 103
$ws.Cells.Item(170, 3).Value = 10
$ws.Cells.Item(170, 4).Value = 0
$ws.Cells.Item(170, 5).Value = "PK"

$ws.Cells.Item(171, 1).Value = 335
$ws.Cells.Item(171, 3).Value = 4
$ws.Cells.Item(171, 4).Value = 4
$ws.Cells.Item(171, 5).Value = "PK"

$ws.Cells.Item(172, 1).Value = 202
$ws.Cells.Item(172, 3).Value = 2
$ws.Cells.Item(172, 4).Value = 0
$ws.Cells.Item(172, 5).Value = "PK"

$ws.Cells.Item(173, 1).Value = 339
$ws.Cells.Item(173, 3).Value = 27
$ws.Cells.Item(173, 4).Value = 3
$ws.Cells.Item(173, 5).Value = "PK"

$ws.Cells.Item(174, 1).Value = 211
$ws.Cells.Item(174, 3).Value = 6
$ws.Cells.Item(174, 4).Value = 0
$ws.Cells.Item(174, 5).Value = "KL"

$ws.Cells.Item(175, 1).Value = 317
$ws.Cells.Item(175, 3).Value = 36
$ws.Cells.Item(175, 4).Value = 0
$ws.Cells.Item(175, 5).Value = "PK"

$ws.Cells.Item(176, 1).Value = 16
$ws.Cells.Item(176, 3).Value = 98
$ws.Cells.Item(176, 4).Value = 12
$ws.Cells.Item(176, 5).Value = "KL"

$ws.Cells.Item(177, 1).Value = 339
$ws.Cells.Item(177, 3).Value = 42
$ws.Cells.Item(177, 4).Value = 1
$ws.Cells.Item(177, 5).Value = "KL"

# --- view state: land on D9, same as the saved workbook -----------------
$ws.Range("D9").Select()
